$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 2.726225666666667
$ws.Range("N2").Value = 8.178677
$ws.Range("O2").Value = 0.0458192930225727
$ws.Range("P2").Value = 0.0458192930225727
$ws.Range("Q2").Value = 97.890527762261
$ws.Range("R2").Value = 881.014749860349
$ws.Range("S2").Value = 0.02798224155389508
$ws.Range("T2").Value = 0.02798224155389508
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.5458848509810237
$ws.Range("P3").Value = 0.5458848509810238
$ws.Range("Q3").Value = 1166.25448877245
$ws.Range("R3").Value = 10496.29039895205
$ws.Range("S3").Value = 0.3333766357599149
$ws.Range("T3").Value = 0.3333766357599149
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("N4").Value = 72.88021500000001
$ws.Range("O4").Value = 0.4082958559964036
$ws.Range("P4").Value = 0.4082958559964037
$ws.Range("Q4").Value = 872.3027831734951
$ws.Range("R4").Value = 7850.725048561455
$ws.Range("S4").Value = 0.2493498374651313
$ws.Range("T4").Value = 0.2493498374651313
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 2.726225666666667
$ws.Range("N5").Value = 8.178677
$ws.Range("O5").Value = 0.0458192930225727
$ws.Range("P5").Value = 0.0458192930225727
$ws.Range("Q5").Value = 46.47839269518179
$ws.Range("R5").Value = 418.3055342566361
$ws.Range("S5").Value = 0.01328595974670767
$ws.Range("T5").Value = 0.01328595974670767
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.5458848509810237
$ws.Range("P6").Value = 0.5458848509810238
$ws.Range("S6").Value = 0.1582871248776892
$ws.Range("T6").Value = 0.1582871248776892
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("N7").Value = 72.88021500000001
$ws.Range("O7").Value = 0.4082958559964036
$ws.Range("P7").Value = 0.4082958559964037
$ws.Range("Q7").Value = 414.1690951335135
$ws.Range("R7").Value = 3727.521856201621
$ws.Range("S7").Value = 0.1183912267010179
$ws.Range("T7").Value = 0.1183912267010179
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 2.726225666666667
$ws.Range("N8").Value = 8.178677
$ws.Range("O8").Value = 0.0458192930225727
$ws.Range("P8").Value = 0.0458192930225727
$ws.Range("Q8").Value = 15.92112517862533
$ws.Range("R8").Value = 143.290126607628
$ws.Range("S8").Value = 0.004551091721969946
$ws.Range("T8").Value = 0.004551091721969946
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.5458848509810237
$ws.Range("P9").Value = 0.5458848509810238
$ws.Range("S9").Value = 0.0542210903434197
$ws.Range("T9").Value = 0.0542210903434197
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("N10").Value = 72.88021500000001
$ws.Range("O10").Value = 0.4082958559964036
$ws.Range("P10").Value = 0.4082958559964037
$ws.Range("S10").Value = 0.04055479183025444
$ws.Range("T10").Value = 0.04055479183025444
